# Update the "想去人数" (want-to-go count) figures in the F column
# for both the "展览" sheet and the "全部类型" sheet, which mirror
# the same underlying data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 575
    $ws.Range("F3").Value = 3617
    $ws.Range("F5").Value = 701
}
